$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill empty homework-score cells (C4:F30) with 0 so every cell has a value
$rng = $ws.Range("C4:F30")
foreach ($cell in $rng.Cells) {
    $v = $cell.Value()
    if ($v -eq $null) {
        $cell.Value = 0
    }
}

# Apply a new style to C4:F30: solid green fill, centered & wrapped text (keep existing thick border)
$rng.Interior.Color = 5296274  # BGR encoding of RGB(146,208,80) == FF92D050
$rng.HorizontalAlignment = -4108  # xlCenter
$rng.VerticalAlignment = -4108    # xlCenter
$rng.WrapText = $true

# Update the active selection on the sheet to L5
$ws.Range("L5").Select()
